# Applies the "aca-metal-tier" StructureDefinition refresh:
#  - Metadata sheet: Version 5.0.0 -> 6.0.0, Date bump, Publisher filled in,
#    the duplicated "Contact / No display for ContactDetail" rows collapsed
#    into a single "Jurisdiction / United States of America" row.
#  - Elements sheet: root element's Short/Definition text updated to describe
#    the ACA Metal Tier extension instead of the generic "Extension" text.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version
$meta.Range("B3").Value = "6.0.0"

# Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher (was blank)
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was the duplicate "Contact" / "No display for ContactDetail" row - remove it entirely,
# shifting every following row up by one.
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

# Root element row (row 2): Short + Definition now describe the extension itself.
$elements.Range("K2").Value = "ACA Metal Tier"
$elements.Range("L2").Value = "Code for Affordable Care Act (ACA) metal tier of the associated plan"
